# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G (header "K", row 1) holds the newly-computed "K" stat for each
# saved-row record (rows 2..64). This replaces the previously-saved
# "Strike#" derived values with the regenerated K values and writes them
# back in place, one per data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K (column G) value, as regenerated by the save routine
$sVals = [ordered]@{
    2  = 4
    3  = 1
    4  = 2
    5  = 1
    6  = 2
    7  = 2
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 0
    13 = 1
    14 = 2
    15 = 2
    16 = 1
    17 = 2
    18 = 2
    19 = 1
    20 = 3
    21 = 1
    22 = 0
    23 = 2
    24 = 1
    25 = 0
    26 = 2
    27 = 0
    28 = 1
    29 = 2
    30 = 1
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 0
    36 = 1
    37 = 1
    38 = 2
    39 = 0
    40 = 1
    41 = 0
    42 = 1
    43 = 1
    44 = 1
    45 = 0
    46 = 1
    47 = 0
    48 = 1
    49 = 1
    50 = 1
    51 = 0
    52 = 1
    53 = 0
    54 = 0
    55 = 1
    56 = 2
    57 = 2
    58 = 2
    59 = 1
    60 = 1
    61 = 0
    62 = 2
    63 = 0
    64 = 0
}

foreach ($row in $sVals.Keys) {
    $ws.Cells.Item($row, 7).Value = $sVals[$row]
}
